$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update duration value
$ws.Range("B2").Value = 10

# Update observation text
$ws.Range("B12").Value = "okk"

# Update data table (rows 16-25, columns A-D)
$data = @(
    @(14, 8, 29, 4),
    @(15, 28, 25, 25),
    @(25, 21, 38, 37),
    @(28, 27, 33, 4),
    @(21, 19, 17, 37),
    @(32, 28, 24, 30),
    @(16, 24, 16, 33),
    @(27, 8, 25, 24),
    @(27, 32, 20, 14),
    @(3, 12, 4, 17)
)

$rowIndex = 16
foreach ($rowValues in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowValues[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowValues[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowValues[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowValues[3]
    $rowIndex++
}
